$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text stays in A1
$ws.Range("A1").Value = "Ingrediente"

# Clear out the old checkbox rows (A2, A5:A11) entirely
$ws.Range("A2").ClearContents()
$ws.Range("A5:A11").ClearContents()

# Ingredients list -> A3
$ingredientes = @'
500 gramos de hongos frescos                                                                                                        
1½ litros de agua o fondo de verduras                                                                                                        
1 pieza de cebolla mediana picada                                                                                                        
2 dientes de ajo picados                                                                                                        
2 cucharadas soperas de crema ácida                                                                                                        
2 cucharadas soperas de mantequilla o margarina                                                                                                        
sal y pimienta al gusto
3 hojas de epazote picado                                                                                                        
1 pizca de tomillo                                                                                                        
2 hojas de laurel
'@
$ws.Range("A3").Value = $ingredientes
$ws.Range("A3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 409.6

# Recipe info block -> A4
$infoBloque = @'
4 comensales
30m
Entrante
Dificultad baja
'@
$ws.Range("A4").Value = $infoBloque
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 57.6

# New empty styled cell at A14
$ws.Range("A14").WrapText = $true

# Update selection
$ws.Range("A3").Select()
